$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 5045.3335
$ws.Range("I70").Value = 2458.7693
$ws.Range("J70").Value = 9248.5
$ws.Range("K70").Value = 7376.3079
$ws.Range("L70").Value = 27745.5
$ws.Range("M70").Value = -7106.3079
$ws.Range("N70").Value = -28285.5
# Row 73
$ws.Range("H73").Value = 5045.3335
$ws.Range("I73").Value = 2458.7693
$ws.Range("J73").Value = 9248.5
$ws.Range("K73").Value = 7376.3079
$ws.Range("L73").Value = 27745.5
$ws.Range("M73").Value = -6440.3079
$ws.Range("N73").Value = -29617.5
# Row 80
$ws.Range("H80").Value = 1940.1111
$ws.Range("I80").Value = 718.8570999999999
$ws.Range("J80").Value = 2717.2727
$ws.Range("K80").Value = 2156.5713
$ws.Range("L80").Value = 8151.8181
$ws.Range("M80").Value = -1158.5713
$ws.Range("N80").Value = -10147.8181
# Row 83
$ws.Range("H83").Value = 1940.1111
$ws.Range("I83").Value = 718.8570999999999
$ws.Range("J83").Value = 2717.2727
$ws.Range("K83").Value = 6469.7139
$ws.Range("L83").Value = 24455.4543
$ws.Range("M83").Value = -1477.7139
$ws.Range("N83").Value = -34439.4543
# Row 112
$ws.Range("H112").Value = 1086.1072
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 1054.2693
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 3162.8079
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -5378.8079
# Row 116
$ws.Range("H116").Value = 11611.6
$ws.Range("I116").Value = 14063.75
$ws.Range("J116").Value = 1803
$ws.Range("K116").Value = 14063.75
$ws.Range("L116").Value = 1803
$ws.Range("M116").Value = -10621.75
$ws.Range("N116").Value = -8687
# Row 127
$ws.Range("H127").Value = 1063.8
$ws.Range("I127").Value = 537.6
$ws.Range("J127").Value = 1590
$ws.Range("K127").Value = 1612.8
$ws.Range("L127").Value = 4770
$ws.Range("M127").Value = 3347.2
$ws.Range("N127").Value = -14690
# Row 129
$ws.Range("H129").Value = 980.89795
$ws.Range("I129").Value = 362
$ws.Range("J129").Value = 1120.15
$ws.Range("K129").Value = 1086
$ws.Range("L129").Value = 3360.45
$ws.Range("M129").Value = 3914
$ws.Range("N129").Value = -13360.45
# Row 131
$ws.Range("H131").Value = 4571.07
$ws.Range("I131").Value = 600.2222
$ws.Range("J131").Value = 4963.791
$ws.Range("K131").Value = 1800.6666
$ws.Range("L131").Value = 14891.373
$ws.Range("M131").Value = 3239.3334
$ws.Range("N131").Value = -24971.373
# Row 138
$ws.Range("H138").Value = 3348.1765
$ws.Range("I138").Value = 1333.6072
$ws.Range("J138").Value = 5800.696
$ws.Range("K138").Value = 4000.8216
$ws.Range("L138").Value = 17402.088
$ws.Range("M138").Value = 1139.1784
$ws.Range("N138").Value = -27682.088

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3656.8289
$ws.Range("I32").Value = 1863.7246
$ws.Range("J32").Value = 21331.715
$ws.Range("K32").Value = 1863.7246
$ws.Range("L32").Value = 21331.715
$ws.Range("M32").Value = -1576.7246
$ws.Range("N32").Value = -21905.715
# Row 61
$ws.Range("H61").Value = 2608.5818
$ws.Range("I61").Value = 1816.7222
$ws.Range("J61").Value = 4108.9473
$ws.Range("K61").Value = 1816.7222
$ws.Range("L61").Value = 4108.9473
$ws.Range("M61").Value = -1604.7222
$ws.Range("N61").Value = -4532.9473
# Row 88
$ws.Range("H88").Value = 7214.5264
$ws.Range("I88").Value = 7487.1665
$ws.Range("J88").Value = 2307
$ws.Range("K88").Value = 7487.1665
$ws.Range("L88").Value = 2307
$ws.Range("M88").Value = -7081.1665
$ws.Range("N88").Value = -3119
# Row 91
$ws.Range("H91").Value = 7214.5264
$ws.Range("I91").Value = 7487.1665
$ws.Range("J91").Value = 2307
$ws.Range("K91").Value = 7487.1665
$ws.Range("L91").Value = 2307
$ws.Range("M91").Value = -6083.1665
$ws.Range("N91").Value = -5115
# Row 110
$ws.Range("H110").Value = 34343.777
$ws.Range("I110").Value = 50680.332
$ws.Range("J110").Value = 1670.6666
$ws.Range("K110").Value = 50680.332
$ws.Range("L110").Value = 1670.6666
$ws.Range("M110").Value = -48635.332
$ws.Range("N110").Value = -5760.6666
# Row 136
$ws.Range("H136").Value = 2608.5818
$ws.Range("I136").Value = 1816.7222
$ws.Range("J136").Value = 4108.9473
$ws.Range("K136").Value = 5450.1666
$ws.Range("L136").Value = 12326.8419
$ws.Range("M136").Value = -2900.1666
$ws.Range("N136").Value = -17426.8419

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2013
$ws.Range("I86").Value = 2018.909
$ws.Range("K86").Value = 2018.909
$ws.Range("M86").Value = -895.9090000000001
# Row 89
$ws.Range("H89").Value = 2013
$ws.Range("I89").Value = 2018.909
$ws.Range("K89").Value = 10094.545
$ws.Range("M89").Value = -4478.545
# Row 94
$ws.Range("H94").Value = 700
$ws.Range("I94").Value = 500
$ws.Range("J94").Value = 900
$ws.Range("K94").Value = 500
$ws.Range("L94").Value = 900
$ws.Range("M94").Value = -49
$ws.Range("N94").Value = -1802
# Row 99
$ws.Range("H99").Value = 1716.3636
$ws.Range("I99").Value = 1570.6666
$ws.Range("J99").Value = 2028.5714
$ws.Range("K99").Value = 1570.6666
$ws.Range("L99").Value = 2028.5714
$ws.Range("M99").Value = -72.66660000000002
$ws.Range("N99").Value = -5024.5714
# Row 134
$ws.Range("H134").Value = 1499.6444
$ws.Range("I134").Value = 1021.9231
$ws.Range("J134").Value = 4604.8335
$ws.Range("K134").Value = 3065.7693
$ws.Range("L134").Value = 13814.5005
$ws.Range("M134").Value = -530.7692999999999
$ws.Range("N134").Value = -18884.5005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 322.22223
$ws.Range("I22").Value = 357.14285
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 357.14285
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -7.14285000000001
$ws.Range("N22").Value = -900
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = $null
# Row 122
$ws.Range("H122").Value = 1094.7391
$ws.Range("I122").Value = 1109.1666
$ws.Range("J122").Value = 1079
$ws.Range("K122").Value = 3327.4998
$ws.Range("L122").Value = 3237
$ws.Range("M122").Value = -877.4998000000001
$ws.Range("N122").Value = -8137

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 8010
$ws.Range("I56").Value = 8010
$ws.Range("K56").Value = 8010
$ws.Range("M56").Value = -7480
# Row 123
$ws.Range("H123").Value = 1892.4445
$ws.Range("I123").Value = 849.75
$ws.Range("J123").Value = 2726.6
$ws.Range("K123").Value = 2549.25
$ws.Range("L123").Value = 8179.799999999999
$ws.Range("M123").Value = -99.25
$ws.Range("N123").Value = -13079.8
# Row 125
$ws.Range("H125").Value = 1321.875
$ws.Range("I125").Value = 710
$ws.Range("J125").Value = 1600
$ws.Range("K125").Value = 2130
$ws.Range("L125").Value = 4800
$ws.Range("M125").Value = 2790
$ws.Range("N125").Value = -14640

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 160.41667
$ws.Range("I2").Value = 50
$ws.Range("J2").Value = 215.625
$ws.Range("K2").Value = 50
$ws.Range("L2").Value = 215.625
$ws.Range("M2").Value = 63
$ws.Range("N2").Value = -441.625
# Row 80
$ws.Range("H80").Value = 3425
$ws.Range("I80").Value = 2858.8635
$ws.Range("J80").Value = 4314.643
$ws.Range("K80").Value = 2858.8635
$ws.Range("L80").Value = 4314.643
$ws.Range("M80").Value = -1860.8635
$ws.Range("N80").Value = -6310.643
# Row 83
$ws.Range("H83").Value = 3425
$ws.Range("I83").Value = 2858.8635
$ws.Range("J83").Value = 4314.643
$ws.Range("K83").Value = 14294.3175
$ws.Range("L83").Value = 21573.215
$ws.Range("M83").Value = -9302.317499999999
$ws.Range("N83").Value = -31557.215

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 342.625
$ws.Range("I55").Value = 100
$ws.Range("K55").Value = 100
$ws.Range("M55").Value = 73
# Row 97
$ws.Range("H97").Value = 62985
$ws.Range("J97").Value = 62985
$ws.Range("L97").Value = 62985
$ws.Range("N97").Value = -64967
# Row 122
$ws.Range("H122").Value = 10103310
$ws.Range("I122").Value = 15874902
$ws.Range("J122").Value = 3024.75
$ws.Range("K122").Value = 47624706
$ws.Range("L122").Value = 9074.25
$ws.Range("M122").Value = -47622256
$ws.Range("N122").Value = -13974.25
# Row 136
$ws.Range("H136").Value = 8548944
$ws.Range("I136").Value = 727.4286
$ws.Range("J136").Value = 18521864
$ws.Range("K136").Value = 2182.2858
$ws.Range("L136").Value = 55565592
$ws.Range("M136").Value = 367.7142000000003
$ws.Range("N136").Value = -55570692

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 2008749.8
$ws.Range("I45").Value = 10000000
$ws.Range("J45").Value = 10937.25
$ws.Range("K45").Value = 10000000
$ws.Range("L45").Value = 10937.25
$ws.Range("M45").Value = -9999509
$ws.Range("N45").Value = -11919.25
# Row 136
$ws.Range("H136").Value = 3056.8235
$ws.Range("I136").Value = 746
$ws.Range("J136").Value = 3767.8462
$ws.Range("K136").Value = 2238
$ws.Range("L136").Value = 11303.5386
$ws.Range("M136").Value = 312
$ws.Range("N136").Value = -16403.5386
